# chore: update Sheets via scheduled runner
# Refreshes cached marketboard price/profit figures (columns H:N) on a
# handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3259.4614
$ws.Range("I132").Value = 3163
$ws.Range("J132").Value = 3476.5
$ws.Range("K132").Value = 9489
$ws.Range("L132").Value = 10429.5
$ws.Range("M132").Value = -6959
$ws.Range("N132").Value = -15489.5
$ws.Range("H137").Value = 2584.875
$ws.Range("I137").Value = 1646.9
$ws.Range("J137").Value = 3254.8572
$ws.Range("K137").Value = 4940.700000000001
$ws.Range("L137").Value = 9764.571599999999
$ws.Range("M137").Value = -2390.700000000001
$ws.Range("N137").Value = -14864.5716
$ws.Range("H138").Value = 3627.7302
$ws.Range("I138").Value = 2040.7037
$ws.Range("J138").Value = 4818
$ws.Range("K138").Value = 6122.1111
$ws.Range("L138").Value = 14454
$ws.Range("M138").Value = -982.1111000000001
$ws.Range("N138").Value = -24734
$ws.Range("H141").Value = 4977.8965
$ws.Range("I141").Value = 1584.25
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 4752.75
$ws.Range("L141").Value = 300000
$ws.Range("M141").Value = 427.25
$ws.Range("N141").Value = -310360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 31325
$ws.Range("J44").Value = 31325
$ws.Range("L44").Value = 31325
$ws.Range("N44").Value = -32301
$ws.Range("H55").Value = 29975
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 29975
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 29975
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -30605
$ws.Range("H61").Value = 2991.3704
$ws.Range("I61").Value = 2960.4285
$ws.Range("J61").Value = 3099.6667
$ws.Range("K61").Value = 2960.4285
$ws.Range("L61").Value = 3099.6667
$ws.Range("M61").Value = -2748.4285
$ws.Range("N61").Value = -3523.6667
$ws.Range("H74").Value = 1560.8223
$ws.Range("I74").Value = 1380.5
$ws.Range("J74").Value = 2118.182
$ws.Range("K74").Value = 1380.5
$ws.Range("L74").Value = 2118.182
$ws.Range("M74").Value = -506.5
$ws.Range("N74").Value = -3866.182
$ws.Range("H77").Value = 1560.8223
$ws.Range("I77").Value = 1380.5
$ws.Range("J77").Value = 2118.182
$ws.Range("K77").Value = 6902.5
$ws.Range("L77").Value = 10590.91
$ws.Range("M77").Value = -2534.5
$ws.Range("N77").Value = -19326.91
$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 5324.485
$ws.Range("I132").Value = 5410.7744
$ws.Range("J132").Value = 3987
$ws.Range("K132").Value = 16232.3232
$ws.Range("L132").Value = 11961
$ws.Range("M132").Value = -13702.3232
$ws.Range("N132").Value = -17021
$ws.Range("H133").Value = 37000
$ws.Range("J133").Value = 37000
$ws.Range("L133").Value = 37000
$ws.Range("N133").Value = -42060
$ws.Range("H136").Value = 2991.3704
$ws.Range("I136").Value = 2960.4285
$ws.Range("J136").Value = 3099.6667
$ws.Range("K136").Value = 8881.2855
$ws.Range("L136").Value = 9299.000100000001
$ws.Range("M136").Value = -6331.2855
$ws.Range("N136").Value = -14399.0001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 17166.666
$ws.Range("I12").Value = 17166.666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 17166.666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -16998.666
$ws.Range("N12").ClearContents()
$ws.Range("H134").Value = 2658.9285
$ws.Range("I134").Value = 2623.7273
$ws.Range("J134").Value = 2788
$ws.Range("K134").Value = 7871.1819
$ws.Range("L134").Value = 8364
$ws.Range("M134").Value = -5336.1819
$ws.Range("N134").Value = -13434
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 21000
$ws.Range("I13").Value = 10000
$ws.Range("J13").Value = 32000
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 32000
$ws.Range("M13").Value = -9861
$ws.Range("N13").Value = -32278
$ws.Range("H58").Value = 1026.4038
$ws.Range("I58").Value = 945.1429000000001
$ws.Range("J58").Value = 1367.7
$ws.Range("K58").Value = 945.1429000000001
$ws.Range("L58").Value = 1367.7
$ws.Range("M58").Value = -742.1429000000001
$ws.Range("N58").Value = -1773.7
$ws.Range("H136").Value = 1026.4038
$ws.Range("I136").Value = 945.1429000000001
$ws.Range("J136").Value = 1367.7
$ws.Range("K136").Value = 2835.4287
$ws.Range("L136").Value = 4103.1
$ws.Range("M136").Value = -285.4287000000004
$ws.Range("N136").Value = -9203.1
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -16068
$ws.Range("H61").Value = 716.7778
$ws.Range("I61").Value = 700.8570999999999
$ws.Range("J61").Value = 772.5
$ws.Range("K61").Value = 2102.5713
$ws.Range("L61").Value = 2317.5
$ws.Range("M61").Value = -1887.5713
$ws.Range("N61").Value = -2747.5
$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 7500
$ws.Range("L80").Value = 22500
$ws.Range("N80").Value = -24372
$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 7500
$ws.Range("L83").Value = 67500
$ws.Range("N83").Value = -76860
$ws.Range("H92").Value = 794
$ws.Range("J92").Value = 794
$ws.Range("L92").Value = 2382
$ws.Range("N92").Value = -4878
$ws.Range("H131").Value = 2924.4717
$ws.Range("J131").Value = 3877.7368
$ws.Range("L131").Value = 11633.2104
$ws.Range("N131").Value = -21713.2104
$ws.Range("H133").Value = 5120
$ws.Range("I133").Value = 2850
$ws.Range("J133").Value = 5532.727
$ws.Range("K133").Value = 8550
$ws.Range("L133").Value = 16598.181
$ws.Range("M133").Value = -3490
$ws.Range("N133").Value = -26718.181
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 28885.5
$ws.Range("J103").Value = 28885.5
$ws.Range("L103").Value = 28885.5
$ws.Range("N103").Value = -31229.5
$ws.Range("H132").Value = 3458
$ws.Range("I132").Value = 3156.889
$ws.Range("K132").Value = 9470.667000000001
$ws.Range("M132").Value = -6940.667000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4463.1055
$ws.Range("I40").Value = 2757
$ws.Range("J40").Value = 5458.3335
$ws.Range("K40").Value = 2757
$ws.Range("L40").Value = 5458.3335
$ws.Range("M40").Value = -2621
$ws.Range("N40").Value = -5730.3335
$ws.Range("H46").Value = 870.7692
$ws.Range("I46").Value = 742.8570999999999
$ws.Range("J46").Value = 1020
$ws.Range("K46").Value = 742.8570999999999
$ws.Range("L46").Value = 1020
$ws.Range("M46").Value = -554.8570999999999
$ws.Range("N46").Value = -1396
$ws.Range("H132").Value = 3205.453
$ws.Range("I132").Value = 3184.9788
$ws.Range("K132").Value = 9554.936399999999
$ws.Range("M132").Value = -7024.936399999999
$ws.Range("H136").Value = 1078.04
$ws.Range("J136").Value = 2625
$ws.Range("L136").Value = 7875
$ws.Range("N136").Value = -12975
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1935.7561
$ws.Range("I132").Value = 1560.877
$ws.Range("J132").Value = 3369.1177
$ws.Range("K132").Value = 4682.630999999999
$ws.Range("L132").Value = 10107.3531
$ws.Range("M132").Value = -2152.630999999999
$ws.Range("N132").Value = -15167.3531
$ws.Range("H136").Value = 1656.9697
$ws.Range("I136").Value = 1518.6538
$ws.Range("J136").Value = 2170.7144
$ws.Range("K136").Value = 4555.9614
$ws.Range("L136").Value = 6512.1432
$ws.Range("M136").Value = -2005.9614
$ws.Range("N136").Value = -11612.1432
